# Updated translation patch to support the huge v2 update.
#
# The sheet was a 4-column (A:D) layout where column A held the source
# text and columns B/C/D occasionally held ad-hoc translation duplicates.
# This collapses it down to a clean 2-column (A:B) layout: column A keeps
# the source text, column B becomes the (currently untranslated, i.e.
# mirrored) translation slot for every row, seeded from whatever
# translation text already existed in B/C/D where present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed column B with column A's value for every used row, so every row
# gets a translation-column entry (most rows had none before).
for ($r = 1; $r -le 94; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $aVal
}

# Rows 4 and 9 already carried a distinct translated value in column C/D
# (not just a duplicate of A) - preserve that actual translation in B
# instead of the mirrored value from the loop above.
$ws.Cells.Item(4, 2).Value = $ws.Cells.Item(4, 3).Value()
$ws.Cells.Item(9, 2).Value = $ws.Cells.Item(9, 4).Value()

# Columns C and D are no longer needed - remove them so the sheet
# collapses to A:B and the dimension shrinks accordingly.
$ws.Columns("C:D").Delete()
